$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.242.59"
$ws.Range("E2").Value = "  +1.71%  "

# Row 3
$ws.Range("D3").Value = "2.612.74"
$ws.Range("E3").Value = "  +1.10%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "561.68"
$ws.Range("E5").Value = "  -1.01%  "

# Row 6
$ws.Range("D6").Value = "142.91"
$ws.Range("E6").Value = "  -0.45%  "

# Row 7
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.22%  "

# Row 8
$ws.Range("D8").Value = "0.602"
$ws.Range("E8").Value = "  +0.30%  "

# Row 9
$ws.Range("D9").Value = "2.637.26"
$ws.Range("E9").Value = "  +1.63%  "

# Row 10
$ws.Range("D10").Value = "6.71"
$ws.Range("E10").Value = "  +0.78%  "

# Row 11
$ws.Range("D11").Value = "0.106"
$ws.Range("E11").Value = "  +1.78%  "

# Row 12
$ws.Range("D12").Value = "0.160"
$ws.Range("E12").Value = "  +3.90%  "

# Row 13
$ws.Range("D13").Value = "0.370"
$ws.Range("E13").Value = "  +7.96%  "

# Row 14
$ws.Range("D14").Value = "3.078.59"
$ws.Range("E14").Value = "  +1.28%  "

# Row 15
$ws.Range("D15").Value = "60.162.54"
$ws.Range("E15").Value = "  +1.49%  "

# Row 16
$ws.Range("D16").Value = "23.49"
$ws.Range("E16").Value = "  +4.71%  "

# Row 17
$ws.Range("D17").Value = "0.0000139"
$ws.Range("E17").Value = "  +1.60%  "

# Row 18
$ws.Range("D18").Value = "2.618.55"
$ws.Range("E18").Value = "  +0.91%  "

# Row 19
$ws.Range("D19").Value = "4.65"
$ws.Range("E19").Value = "  +3.05%  "

# Row 20
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "10.85"
$ws.Range("E20").Value = "  +6.32%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "344.67"
$ws.Range("E21").Value = "  +2.57%  "

# Row 22
$ws.Range("D22").Value = "6.91"
$ws.Range("E22").Value = "  +11.58%  "

# Row 23
$ws.Range("E23").Value = "  +0.30%  "

# Row 24
$ws.Range("D24").Value = "0.523"
$ws.Range("E24").Value = "  +16.15%  "

# Row 25
$ws.Range("D25").Value = "62.80"
$ws.Range("E25").Value = "  -2.09%  "

# Row 26
$ws.Range("D26").Value = "0.994"
$ws.Range("E26").Value = "  -0.40%  "

# Row 27
$ws.Range("D27").Value = "0.161"
$ws.Range("E27").Value = "  -0.49%  "

# Row 28
$ws.Range("D28").Value = "7.71"
$ws.Range("E28").Value = "  +6.30%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0792"
$ws.Range("E29").Value = "  +0.98%  "

# Row 30
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.06%  "

# Row 31
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "6.28"
$ws.Range("E31").Value = "  +3.84%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "1.70"
$ws.Range("E32").Value = "  +1.38%  "

# Row 33
$ws.Range("D33").Value = "160.05"
$ws.Range("E33").Value = "  +1.45%  "

# Row 34
$ws.Range("D34").Value = "19.51"
$ws.Range("E34").Value = "  +2.55%  "

# Row 35
$ws.Range("D35").Value = "4.23"
$ws.Range("E35").Value = "  +4.81%  "

# Row 36
$ws.Range("D36").Value = "0.954"
$ws.Range("E36").Value = "  +7.63%  "

# Row 37
$ws.Range("D37").Value = "1.20"
$ws.Range("E37").Value = "  +5.24%  "

# Row 38
$ws.Range("D38").Value = "1.56"
$ws.Range("E38").Value = "  +3.59%  "

# Row 39
$ws.Range("D39").Value = "37.78"
$ws.Range("E39").Value = "  +2.44%  "

# Row 40
$ws.Range("D40").Value = "0.859"
$ws.Range("E40").Value = "  -2.31%  "

# Row 41
$ws.Range("D41").Value = "3.76"
$ws.Range("E41").Value = "  +3.11%  "

# Row 42
$ws.Range("D42").Value = "300.62"
$ws.Range("E42").Value = "  +2.18%  "

# Row 43
$ws.Range("D43").Value = "141.42"
$ws.Range("E43").Value = "  +13.69%  "

# Row 44
$ws.Range("D44").Value = "0.995"
$ws.Range("E44").Value = "  -0.34%  "

# Row 45
$ws.Range("D45").Value = "0.0982"
$ws.Range("E45").Value = "  +0.80%  "

# Row 46
$ws.Range("D46").Value = "0.605"
$ws.Range("E46").Value = "  +1.21%  "

# Row 47
$ws.Range("D47").Value = "0.0242"
$ws.Range("E47").Value = "  +4.38%  "

# Row 48
$ws.Range("D48").Value = "0.0545"
$ws.Range("E48").Value = "  +1.76%  "

# Row 49
$ws.Range("D49").Value = "10.66"
$ws.Range("E49").Value = "  +0.41%  "

# Row 50
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "4.81"
$ws.Range("E50").Value = "  +6.70%  "

# Row 51
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "19.45"
$ws.Range("E51").Value = "  +4.95%  "
